# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.265.44"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.57%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.966.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.86%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "44.81"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.03%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06894"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.35"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "109.17"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.29%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07771"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.934.57"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.490"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.7188"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "287.67"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "31.152.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.30%  "

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.33"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.57%  "

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007815"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.42%  "

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.186.19"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.549"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.29%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.600"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.935"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.84"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.09%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.37"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.54%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.227"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.14%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.439"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.586"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.652"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.488"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.04998"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7673"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.194"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02058"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.733"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.01%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.710"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +7.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.449"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +7.82%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4561"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.88"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8862"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "72.58"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.149"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +9.73%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.20%  "

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.492"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.37%  "

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "969.65"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.93%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1273"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.98%  "
